$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 649 — shifts rows 649:690 down to 650:691
# and the sheet dimension grows from A1:D690 to A1:D691 automatically.
$ws.Rows("649:649").Insert()

# Fill the new row 649 with the new reading.
# Column A holds a date-shaped string ("2026/01/16") that must stay literal
# text (matching every other row in the sheet) instead of being
# auto-converted to a date serial by Excel's input parser, so format the
# cell as Text first, then strip the formatting override back off again
# once the literal string is safely stored.
$ws.Range("A649").NumberFormat = "@"
$ws.Range("A649").Value = "2026/01/16"
$ws.Range("A649").ClearFormats()

$ws.Range("B649").Value = "金"
$ws.Range("C649").Value = 14
$ws.Range("D649").Value = 201
